$d = $word.ActiveDocument

# Paragraph 6 = body under "Objetivos" -> numbered list (was under "Programa resumido")
$d.Paragraphs.Item(6).Range.Text = "1.Petróleo: histórico, constituinte, composição e classificação `v2.Geologia do petróleo: origem.`v3.Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; `v4.Perfuração: equipamentos, operações;`v5.Completação e reservatórios: tipos, etapas, equipamentos, reservatórios;`v6.Elevação: elevação natural, sistemas de bombeamentos;`v7.Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural;`v8.Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais;`v9.Processos de refino: objetivo, tipos de processos, esquemas de refino.`v9.1- Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo;`v9.2- Desasfaltação: carga, descrição e variáveis do processo, produtos;`v9.3- Coqueamento retardado: carga, descrição e variáveis do processo, produtos;`v9.4- Craqueamento catalítico: carga, descrição e variáveis do processo, produtos;`v9.5- Hidrorrefino: carga, descrição e variáveis do processo, produtos;`v9.6- Reforma catalítica: carga, descrição e variáveis do processo, produtos;`v9.7- Alquilação e isomerização: carga, descrição e variáveis do processo, produtos;`v9.8- Tratamento de derivados: tratamento com aminas, tratamentos cáusticos;`v9.9- Geração de hidrogênio: carga, descrição e variáveis do processo;`v9.10- Recuperação de Enxofre: Processo Claus.`v10.Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo."

# Paragraph 8 (ListBullet under "Docente(s) Responsavel(eis)") -> "Visao integrada..." (was under "Objetivos")
$d.Paragraphs.Item(8).Range.Text = "Visão integrada sobre petróleo e gás natural, desde a origem até o processamento primário. Descrições, características e aplicações dos derivados do petróleo. Processo e esquemas de refino e processamento do gás natural."

# Paragraph 10 (body under "Programa resumido") -> long unnumbered program text (was under "Programa")
$d.Paragraphs.Item(10).Range.Text = "Petróleo: histórico, constituinte, composição e classificação; Geologia do petróleo: origem; Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; Perfuração: equipamentos, operações; Completação e reservatórios: tipos, etapas, equipamentos, reservatórios; Elevação: elevação natural, bombeios; Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural; Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais; Processos de refino: objetivo, tipos de processos, esquemas de refino; Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo; Desasfaltação: carga, descrição e variáveis do processo, produtos; Coqueamento retardado: carga, descrição e variáveis do processo, produtos; Craqueamento catalítico: carga, descrição e variáveis do processo, produtos; Hidrorrefino: carga, descrição e variáveis do processo, produtos; Reforma catalítica: carga, descrição e variáveis do processo, produtos; Alquilação e isomerização: carga, descrição e variáveis do processo, produtos; Tratamento de derivados: tratamento com aminas, tratamentos cáusticos; Geração de hidrogênio: carga, descrição e variáveis do processo; Recuperação de Enxofre: Processo Claus; Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo."

# Paragraph 12 (body under "Programa") -> "Aulas expositivas..." (was the Metodo value under Avaliacao)
$d.Paragraphs.Item(12).Range.Text = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários"

# Paragraph 14 (Avaliacao ListBullet): the Metodo/Criterio/Norma values rotate forward
# (Criterio->Metodo, Norma->Criterio, Bibliografia list->Norma). Replace starting from the
# LAST value first so an old value is never re-matched after being freshly (re)written, and
# re-fetch Paragraphs.Item(14).Range fresh each time since Find/Replace collapses the Range
# used to the replacement span, not the whole paragraph.
$d.Paragraphs.Item(14).Range.Find.Execute("Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.", $false, $false, $false, $false, $false, $true, 1, $false, "a)Speight, J. G., The Chemistry and Technology of Petroleum, CRC Press, 4ª Edição, 2007;`vb)Thomas, J. E. (Organizador), Fundamentos de Engenharia de Petróleo, Editora Interciência, 2ª Edição, 2004;`vc)Brasil, N. I., Araújo, M. A. S., Souza, E. C. M, Processamento de Petróleo e Gás, Editora LTC, 1ª Edição, 2012;`vd)Fundamentos do Refino do Petróleo  Tecnologia e Economia, Szklo, A. S., Uller, V. C., Bonfá, M. H. P., Editora Interciência, 3ª Edição, 2012.`ve)Oil and Gas Journal;`vf)Revista Petro & Química.", 2) | Out-Null
$d.Paragraphs.Item(14).Range.Find.Execute("Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.", $false, $false, $false, $false, $false, $true, 1, $false, "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.", 2) | Out-Null
$d.Paragraphs.Item(14).Range.Find.Execute("Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários", $false, $false, $false, $false, $false, $true, 1, $false, "Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.", 2) | Out-Null

# Paragraph 16 (body under "Bibliografia") -> "1285870 - Marcos Villela Barcza" (was under "Docente(s)")
$d.Paragraphs.Item(16).Range.Text = "1285870 - Marcos Villela Barcza"

Write-Output "done"
